$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C and E values for rows 2-6 (E2, C3:E3, C4:E4, C5:E5, C6:E6)
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Update C and E values for rows 7-19
$ws.Range("C7").Value = 0.4454453461194552
$ws.Range("E7").Value = 2.887668087172179

$ws.Range("C8").Value = 5.461771395837989
$ws.Range("E8").Value = 3.667647087004666

$ws.Range("C9").Value = 3.857137494611718
$ws.Range("E9").Value = 3.60208822706134

$ws.Range("C10").Value = 4.325618632128836
$ws.Range("E10").Value = 3.998755562728684

$ws.Range("C11").Value = 4.439146757103352
$ws.Range("E11").Value = 3.861679870292711

$ws.Range("C12").Value = 5.723509166364238
$ws.Range("E12").Value = 4.058053416301188

$ws.Range("C13").Value = 4.674926984813466
$ws.Range("E13").Value = 4.726969153629335

$ws.Range("C14").Value = 0.4167846160013644
$ws.Range("E14").Value = 1.125570778878981

$ws.Range("C15").Value = -1.488064879190421
$ws.Range("E15").Value = 3.63609986063671

$ws.Range("C16").Value = 1.712986619197032
$ws.Range("E16").Value = 2.431967849366434

$ws.Range("C17").Value = -0.5717743519535134
$ws.Range("E17").Value = 1.732880403074311

$ws.Range("C18").Value = -0.2867681914691111
$ws.Range("E18").Value = 2.013081730696564

$ws.Range("C19").Value = 0.4432539413513181
$ws.Range("E19").Value = 1.364481450639365
